$d = $word.ActiveDocument

$pairs = @(
    @("2024-03-15 Friday", "2024-03-16 Saturday"),
    @("390÷3=130, 0", "976÷2=488, 0"),
    @("815÷7=116, 3", "308÷4=77, 0"),
    @("873÷7=124, 5", "535÷7=76, 3"),
    @("935÷5=187, 0", "358÷8=44, 6"),
    @("176÷8=22, 0", "604÷7=86, 2"),
    @("172÷2=86, 0", "605÷2=302, 1"),
    @("743÷7=106, 1", "547÷4=136, 3"),
    @("769÷7=109, 6", "129÷9=14, 3"),
    @("970÷8=121, 2", "448÷2=224, 0"),
    @("389÷3=129, 2", "238÷9=26, 4"),
    @("163÷6=27, 1", "251÷5=50, 1"),
    @("596÷5=119, 1", "584÷8=73, 0"),
    @("324÷3=108, 0", "839÷8=104, 7"),
    @("811÷5=162, 1", "530÷2=265, 0"),
    @("330÷2=165, 0", "491÷8=61, 3"),
    @("650÷3=216, 2", "448÷8=56, 0"),
    @("525÷4=131, 1", "597÷8=74, 5"),
    @("483÷8=60, 3", "485÷5=97, 0"),
    @("487÷7=69, 4", "333÷8=41, 5"),
    @("896÷8=112, 0", "286÷2=143, 0"),
    @("797÷5=159, 2", "150÷9=16, 6"),
    @("733÷4=183, 1", "579÷8=72, 3"),
    @("219÷3=73, 0", "908÷7=129, 5"),
    @("843÷8=105, 3", "106÷7=15, 1"),
    @("435÷7=62, 1", "253÷7=36, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
